$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the existing row-1 values (read via Value2 -- Value's getter is unreliable
# in this host) before they get overwritten, since columns B..I shift right to D..K.
$oldB1 = $ws.Range("B1").Value2
$oldC1 = $ws.Range("C1").Value2
$oldD1 = $ws.Range("D1").Value2
$oldE1 = $ws.Range("E1").Value2
$oldF1 = $ws.Range("F1").Value2
$oldG1 = $ws.Range("G1").Value2
$oldH1 = $ws.Range("H1").Value2
$oldI1 = $ws.Range("I1").Value2

# New leading columns: player's last name, first name, position.
$ws.Range("A1").Value = "Walford"
$ws.Range("B1").Value = "Clive"
$ws.Range("C1").Value = "TE"

# The rest of the original row slides three columns to the right (D..K).
# Force text format on the numeric-looking values so they remain strings
# ("2018-12-30", "16", "27.070") instead of being coerced to date/number types.
$ws.Range("D1").NumberFormat = "@"
$ws.Range("E1").NumberFormat = "@"
$ws.Range("F1").NumberFormat = "@"

$ws.Range("D1").Value = $oldB1
$ws.Range("E1").Value = $oldC1
$ws.Range("F1").Value = $oldD1
$ws.Range("G1").Value = $oldE1
$ws.Range("H1").Value = $oldF1
$ws.Range("I1").Value = $oldG1
$ws.Range("J1").Value = $oldH1
$ws.Range("K1").Value = $oldI1

# New trailing numeric column.
$ws.Range("L1").Value = 0
